$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$t = @(0, 0.2105263157894737, 0.4210526315789473, 0.631578947368421, 0.8421052631578947, 1.052631578947368, 1.263157894736842, 1.473684210526316, 1.684210526315789, 1.894736842105263, 2.105263157894737, 2.315789473684211, 2.526315789473684, 2.736842105263158, 2.947368421052631, 3.157894736842105, 3.368421052631579, 3.578947368421052, 3.789473684210526, 4)

$x = @(1, 1.376731301939058, 1.692444042019321, 1.93314541181618, 2.088167206574094, 2.150638654226231, 2.11779093714812, 1.991079909836297, 1.776121573501259, 1.482443167426553, 1.123060909055933, 0.7139030979293716, 0.2731041522685442, -0.1797991325509998, -0.6247334807477848, -1.041978810407272, -1.413042253732919, -1.721477785535777, -1.953615132938157, -2.09916566004411)

for ($i = 0; $i -lt $t.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $t[$i]
    $ws.Cells.Item($row, 2).Value = $x[$i]
}
